$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value2 = 8.644702
$ws.Range("H2").Value2 = 25.934106
$ws.Range("I2").Value2 = 0.2659466972557785
$ws.Range("J2").Value2 = 0.2659466972557785
$ws.Range("M2").Value2 = 9.886733333333334
$ws.Range("N2").Value2 = 29.6602
$ws.Range("O2").Value2 = 0.2800365922084035
$ws.Range("P2").Value2 = 0.2800365922084035
$ws.Range("Q2").Value2 = 85.46786342013334
$ws.Range("R2").Value2 = 769.2107707812
$ws.Range("S2").Value2 = 0.07447480680858819
$ws.Range("T2").Value2 = 0.07447480680858819
$ws.Range("G3").Value2 = 8.644702
$ws.Range("H3").Value2 = 25.934106
$ws.Range("I3").Value2 = 0.2659466972557785
$ws.Range("J3").Value2 = 0.2659466972557785
$ws.Range("M3").Value2 = 9.340016666666665
$ws.Range("O3").Value2 = 0.2645511262738982
$ws.Range("P3").Value2 = 0.2645511262738982
$ws.Range("Q3").Value2 = 80.74166075836666
$ws.Range("R3").Value2 = 726.6749468252999
$ws.Range("S3").Value2 = 0.07035649828783963
$ws.Range("T3").Value2 = 0.07035649828783963
$ws.Range("G4").Value2 = 8.644702
$ws.Range("H4").Value2 = 25.934106
$ws.Range("I4").Value2 = 0.2659466972557785
$ws.Range("J4").Value2 = 0.2659466972557785
$ws.Range("M4").Value2 = 9.811931666666666
$ws.Range("N4").Value2 = 29.435795
$ws.Range("O4").Value2 = 0.277917873808847
$ws.Range("P4").Value2 = 0.277917873808847
$ws.Range("Q4").Value2 = 84.82122530269667
$ws.Range("R4").Value2 = 763.39102772427
$ws.Range("S4").Value2 = 0.07391134064781107
$ws.Range("T4").Value2 = 0.07391134064781107
$ws.Range("G5").Value2 = 8.644702
$ws.Range("H5").Value2 = 25.934106
$ws.Range("I5").Value2 = 0.2659466972557785
$ws.Range("J5").Value2 = 0.2659466972557785
$ws.Range("M5").Value2 = 6.266466333333334
$ws.Range("N5").Value2 = 18.799399
$ws.Range("O5").Value2 = 0.1774944077088512
$ws.Range("P5").Value2 = 0.1774944077088512
$ws.Range("Q5").Value2 = 54.17173404469934
$ws.Range("R5").Value2 = 487.545606402294
$ws.Range("S5").Value2 = 0.04720405151153956
$ws.Range("T5").Value2 = 0.04720405151153956
$ws.Range("I6").Value2 = 0.331579210752513
$ws.Range("J6").Value2 = 0.3315792107525131
$ws.Range("M6").Value2 = 9.886733333333334
$ws.Range("N6").Value2 = 29.6602
$ws.Range("O6").Value2 = 0.2800365922084035
$ws.Range("P6").Value2 = 0.2800365922084035
$ws.Range("Q6").Value2 = 106.5603257719556
$ws.Range("R6").Value2 = 959.0429319476002
$ws.Range("S6").Value2 = 0.09285431222628578
$ws.Range("T6").Value2 = 0.0928543122262858
$ws.Range("I7").Value2 = 0.331579210752513
$ws.Range("J7").Value2 = 0.3315792107525131
$ws.Range("M7").Value2 = 9.340016666666665
$ws.Range("O7").Value2 = 0.2645511262738982
$ws.Range("P7").Value2 = 0.2645511262738982
$ws.Range("R7").Value2 = 906.0097674769
$ws.Range("S7").Value2 = 0.08771965365358758
$ws.Range("T7").Value2 = 0.0877196536535876
$ws.Range("I8").Value2 = 0.331579210752513
$ws.Range("J8").Value2 = 0.3315792107525131
$ws.Range("M8").Value2 = 9.811931666666666
$ws.Range("N8").Value2 = 29.435795
$ws.Range("O8").Value2 = 0.277917873808847
$ws.Range("P8").Value2 = 0.277917873808847
$ws.Range("Q8").Value2 = 105.7541049809678
$ws.Range("R8").Value2 = 951.78694482871
$ws.Range("S8").Value2 = 0.09215178925155398
$ws.Range("T8").Value2 = 0.092151789251554
$ws.Range("I9").Value2 = 0.331579210752513
$ws.Range("J9").Value2 = 0.3315792107525131
$ws.Range("M9").Value2 = 6.266466333333334
$ws.Range("N9").Value2 = 18.799399
$ws.Range("O9").Value2 = 0.1774944077088512
$ws.Range("P9").Value2 = 0.1774944077088512
$ws.Range("Q9").Value2 = 67.54068016254023
$ws.Range("R9").Value2 = 607.8661214628621
$ws.Range("S9").Value2 = 0.05885345562108565
$ws.Range("T9").Value2 = 0.05885345562108566
$ws.Range("G10").Value2 = 7.273908666666667
$ws.Range("H10").Value2 = 21.821726
$ws.Range("I10").Value2 = 0.2237754391117454
$ws.Range("J10").Value2 = 0.2237754391117454
$ws.Range("M10").Value2 = 9.886733333333334
$ws.Range("N10").Value2 = 29.6602
$ws.Range("O10").Value2 = 0.2800365922084035
$ws.Range("P10").Value2 = 0.2800365922084035
$ws.Range("Q10").Value2 = 71.91519527835557
$ws.Range("R10").Value2 = 647.2367575052001
$ws.Range("S10").Value2 = 0.06266531138879228
$ws.Range("T10").Value2 = 0.06266531138879228
$ws.Range("G11").Value2 = 7.273908666666667
$ws.Range("H11").Value2 = 21.821726
$ws.Range("I11").Value2 = 0.2237754391117454
$ws.Range("J11").Value2 = 0.2237754391117454
$ws.Range("M11").Value2 = 9.340016666666665
$ws.Range("O11").Value2 = 0.2645511262738982
$ws.Range("P11").Value2 = 0.2645511262738982
$ws.Range("Q11").Value2 = 67.93842817847778
$ws.Range("R11").Value2 = 611.4458536063
$ws.Range("S11").Value2 = 0.05920004444944837
$ws.Range("T11").Value2 = 0.05920004444944837
$ws.Range("G12").Value2 = 7.273908666666667
$ws.Range("H12").Value2 = 21.821726
$ws.Range("I12").Value2 = 0.2237754391117454
$ws.Range("J12").Value2 = 0.2237754391117454
$ws.Range("M12").Value2 = 9.811931666666666
$ws.Range("N12").Value2 = 29.435795
$ws.Range("O12").Value2 = 0.277917873808847
$ws.Range("P12").Value2 = 0.277917873808847
$ws.Range("Q12").Value2 = 71.37109478690778
$ws.Range("R12").Value2 = 642.3398530821701
$ws.Range("S12").Value2 = 0.06219119424857737
$ws.Range("T12").Value2 = 0.06219119424857737
$ws.Range("G13").Value2 = 7.273908666666667
$ws.Range("H13").Value2 = 21.821726
$ws.Range("I13").Value2 = 0.2237754391117454
$ws.Range("J13").Value2 = 0.2237754391117454
$ws.Range("M13").Value2 = 6.266466333333334
$ws.Range("N13").Value2 = 18.799399
$ws.Range("O13").Value2 = 0.1774944077088512
$ws.Range("P13").Value2 = 0.1774944077088512
$ws.Range("Q13").Value2 = 45.58170377140823
$ws.Range("R13").Value2 = 410.2353339426741
$ws.Range("S13").Value2 = 0.03971888902492735
$ws.Range("T13").Value2 = 0.03971888902492735
$ws.Range("G14").Value2 = 5.80867
$ws.Range("H14").Value2 = 17.42601
$ws.Range("I14").Value2 = 0.1786986528799631
$ws.Range("J14").Value2 = 0.1786986528799631
$ws.Range("M14").Value2 = 9.886733333333334
$ws.Range("N14").Value2 = 29.6602
$ws.Range("O14").Value2 = 0.2800365922084035
$ws.Range("P14").Value2 = 0.2800365922084035
$ws.Range("Q14").Value2 = 57.42877131133334
$ws.Range("R14").Value2 = 516.8589418020001
$ws.Range("S14").Value2 = 0.05004216178473728
$ws.Range("T14").Value2 = 0.05004216178473729
$ws.Range("G15").Value2 = 5.80867
$ws.Range("H15").Value2 = 17.42601
$ws.Range("I15").Value2 = 0.1786986528799631
$ws.Range("J15").Value2 = 0.1786986528799631
$ws.Range("M15").Value2 = 9.340016666666665
$ws.Range("O15").Value2 = 0.2645511262738982
$ws.Range("P15").Value2 = 0.2645511262738982
$ws.Range("Q15").Value2 = 54.25307461116666
$ws.Range("R15").Value2 = 488.2776715005
$ws.Range("S15").Value2 = 0.04727492988302262
$ws.Range("T15").Value2 = 0.04727492988302263
$ws.Range("G16").Value2 = 5.80867
$ws.Range("H16").Value2 = 17.42601
$ws.Range("I16").Value2 = 0.1786986528799631
$ws.Range("J16").Value2 = 0.1786986528799631
$ws.Range("M16").Value2 = 9.811931666666666
$ws.Range("N16").Value2 = 29.435795
$ws.Range("O16").Value2 = 0.277917873808847
$ws.Range("P16").Value2 = 0.277917873808847
$ws.Range("Q16").Value2 = 56.99427311421667
$ws.Range("R16").Value2 = 512.94845802795
$ws.Range("S16").Value2 = 0.04966354966090453
$ws.Range("T16").Value2 = 0.04966354966090454
$ws.Range("G17").Value2 = 5.80867
$ws.Range("H17").Value2 = 17.42601
$ws.Range("I17").Value2 = 0.1786986528799631
$ws.Range("J17").Value2 = 0.1786986528799631
$ws.Range("M17").Value2 = 6.266466333333334
$ws.Range("N17").Value2 = 18.799399
$ws.Range("O17").Value2 = 0.1774944077088512
$ws.Range("P17").Value2 = 0.1774944077088512
$ws.Range("Q17").Value2 = 36.39983499644334
$ws.Range("R17").Value2 = 327.5985149679901
$ws.Range("S17").Value2 = 0.03171801155129865
$ws.Range("T17").Value2 = 0.03171801155129866
